$p = $ppt.ActivePresentation

# Slide 13
$s = $p.Slides.Add(13, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("were")
$title.InsertAfter(" ")
$title.InsertAfter("the")
$title.InsertAfter(" ")
$title.InsertAfter("courses")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("taken")
$title.InsertAfter(" ")
$title.InsertAfter("that")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("been")
$title.InsertAfter(" ")
$title.InsertAfter("the")
$title.InsertAfter(" ")
$title.InsertAfter("most")
$title.InsertAfter(" ")
$title.InsertAfter("beneficial")
$title.InsertAfter(" ")
$title.InsertAfter("for")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("job?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Applied courses outside Department of Statistics")
$body.InsertAfter("`r")
$body.InsertAfter("Biostatistics, Education")
$body.InsertAfter("`r")
$body.InsertAfter("Very practical advice")
$body.InsertAfter("`r")
$body.InsertAfter("Specialized perspectives ()")
$body.InsertAfter("`r")
$body.InsertAfter("Theory of Mathematical Statistics")
$body.InsertAfter("`r")
$body.InsertAfter("KNowing foundations increases your confidence level")
$body.InsertAfter("`r")
$body.InsertAfter("Also makes you a quick learner of new methods")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 3
$body.Paragraphs(3).IndentLevel = 3
$body.Paragraphs(4).IndentLevel = 3
$body.Paragraphs(5).IndentLevel = 2
$body.Paragraphs(6).IndentLevel = 3
$body.Paragraphs(7).IndentLevel = 3

# Slide 14
$s = $p.Slides.Add(14, 6)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("Which")
$title.InsertAfter(" ")
$title.InsertAfter("aspects")
$title.InsertAfter(" ")
$title.InsertAfter("of")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("work")
$title.InsertAfter(" ")
$title.InsertAfter("did")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("learn")
$title.InsertAfter(" ")
$title.InsertAfter("on")
$title.InsertAfter(" ")
$title.InsertAfter("the")
$title.InsertAfter(" ")
$title.InsertAfter("job?")

# Slide 15
$s = $p.Slides.Add(15, 6)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What,")
$title.InsertAfter(" ")
$title.InsertAfter("if")
$title.InsertAfter(" ")
$title.InsertAfter("anything,")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("wish")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("knew")
$title.InsertAfter(" ")
$title.InsertAfter("before")
$title.InsertAfter(" ")
$title.InsertAfter("entering")
$title.InsertAfter(" ")
$title.InsertAfter("the")
$title.InsertAfter(" ")
$title.InsertAfter("workforce")
$title.InsertAfter(" ")
$title.InsertAfter("or")
$title.InsertAfter(" ")
$title.InsertAfter("taking")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("current")
$title.InsertAfter(" ")
$title.InsertAfter("position?")

# Slide 16
$s = $p.Slides.Add(16, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("How")
$title.InsertAfter(" ")
$title.InsertAfter("much")
$title.InsertAfter(" ")
$title.InsertAfter("should")
$title.InsertAfter(" ")
$title.InsertAfter("an")
$title.InsertAfter(" ")
$title.InsertAfter("entry")
$title.InsertAfter(" ")
$title.InsertAfter("level")
$title.InsertAfter(" ")
$title.InsertAfter("statistician")
$title.InsertAfter(" ")
$title.InsertAfter("expect")
$title.InsertAfter(" ")
$title.InsertAfter("to")
$title.InsertAfter(" ")
$title.InsertAfter("make")
$title.InsertAfter(" ")
$title.InsertAfter("in")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("field?")
$title.InsertAfter(" ")
$title.InsertAfter("After")
$title.InsertAfter(" ")
$title.InsertAfter("5")
$title.InsertAfter(" ")
$title.InsertAfter("years")
$title.InsertAfter(" ")
$title.InsertAfter("of")
$title.InsertAfter(" ")
$title.InsertAfter("experience?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Outside my area of expertise (sorry!)")
$body.InsertAfter("`r")
$body.InsertAfter("See various salary surveys")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 3

# Slide 17
$s = $p.Slides.Add(17, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("recommend")
$title.InsertAfter(" ")
$title.InsertAfter("as")
$title.InsertAfter(" ")
$title.InsertAfter("the")
$title.InsertAfter(" ")
$title.InsertAfter("best")
$title.InsertAfter(" ")
$title.InsertAfter("places")
$title.InsertAfter(" ")
$title.InsertAfter("to")
$title.InsertAfter(" ")
$title.InsertAfter("look")
$title.InsertAfter(" ")
$title.InsertAfter("for")
$title.InsertAfter(" ")
$title.InsertAfter("new")
$title.InsertAfter(" ")
$title.InsertAfter("job")
$title.InsertAfter(" ")
$title.InsertAfter("opportunities?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Outside my area of expertise (sorry!)")
$body.Paragraphs(1).IndentLevel = 2

# Slide 18
$s = $p.Slides.Add(18, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("How")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("find")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("clients/collaborators?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Word of mouth")
$body.InsertAfter("`r")
$body.InsertAfter("Let everyone know you are looking")
$body.InsertAfter("`r")
$body.InsertAfter("Increase your visibility")
$body.InsertAfter("`r")
$body.InsertAfter("Blogs")
$body.InsertAfter("`r")
$body.InsertAfter("Social media")
$body.InsertAfter("`r")
$body.InsertAfter("Presentations")
$body.InsertAfter("`r")
$body.InsertAfter("Volunteer")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 3
$body.Paragraphs(3).IndentLevel = 2
$body.Paragraphs(4).IndentLevel = 3
$body.Paragraphs(5).IndentLevel = 3
$body.Paragraphs(6).IndentLevel = 3
$body.Paragraphs(7).IndentLevel = 3

# Slide 19
$s = $p.Slides.Add(19, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("advice")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("to")
$title.InsertAfter(" ")
$title.InsertAfter("give")
$title.InsertAfter(" ")
$title.InsertAfter("regarding")
$title.InsertAfter(" ")
$title.InsertAfter("balancing")
$title.InsertAfter(" ")
$title.InsertAfter("competing")
$title.InsertAfter(" ")
$title.InsertAfter("priorities?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Life is short. Do the fun stuff first.")
$body.InsertAfter("`r")
$body.InsertAfter("Ask your boss for priorities")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 2

# Slide 20
$s = $p.Slides.Add(20, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("when")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("a")
$title.InsertAfter(" ")
$title.InsertAfter("project")
$title.InsertAfter(" ")
$title.InsertAfter("that")
$title.InsertAfter(" ")
$title.InsertAfter("is")
$title.InsertAfter(" ")
$title.InsertAfter("outside")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("level")
$title.InsertAfter(" ")
$title.InsertAfter("of")
$title.InsertAfter(" ")
$title.InsertAfter("expertise")
$title.InsertAfter(" ")
$title.InsertAfter("(who")
$title.InsertAfter(" ")
$title.InsertAfter("or")
$title.InsertAfter(" ")
$title.InsertAfter("what")
$title.InsertAfter(" ")
$title.InsertAfter("resources")
$title.InsertAfter(" ")
$title.InsertAfter("do")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("consult)?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Nothing is ever outside my area of expertise!")
$body.InsertAfter("`r")
$body.InsertAfter("Just kidding")
$body.InsertAfter("`r")
$body.InsertAfter("But I do try to know a bit about everything")
$body.InsertAfter("`r")
$body.InsertAfter("Opportunity to learn on the job")
$body.InsertAfter("`r")
$body.InsertAfter("Books, short courses, Internet resources")
$body.InsertAfter("`r")
$body.InsertAfter("Build a network of helpers")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 3
$body.Paragraphs(3).IndentLevel = 3
$body.Paragraphs(4).IndentLevel = 2
$body.Paragraphs(5).IndentLevel = 3
$body.Paragraphs(6).IndentLevel = 3

# Slide 21
$s = $p.Slides.Add(21, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("is")
$title.InsertAfter(" ")
$title.InsertAfter("an")
$title.InsertAfter(" ")
$title.InsertAfter("issue")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("faced")
$title.InsertAfter(" ")
$title.InsertAfter("working")
$title.InsertAfter(" ")
$title.InsertAfter("with")
$title.InsertAfter(" ")
$title.InsertAfter("a")
$title.InsertAfter(" ")
$title.InsertAfter("client")
$title.InsertAfter(" ")
$title.InsertAfter("or")
$title.InsertAfter(" ")
$title.InsertAfter("collaborator,")
$title.InsertAfter(" ")
$title.InsertAfter("and")
$title.InsertAfter(" ")
$title.InsertAfter("what")
$title.InsertAfter(" ")
$title.InsertAfter("did")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("learn")
$title.InsertAfter(" ")
$title.InsertAfter("from")
$title.InsertAfter(" ")
$title.InsertAfter("it?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Get and give regular feedback")
$body.Paragraphs(1).IndentLevel = 2

# Slide 22
$s = $p.Slides.Add(22, 2)
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.InsertAfter("What")
$title.InsertAfter(" ")
$title.InsertAfter("are")
$title.InsertAfter(" ")
$title.InsertAfter("some")
$title.InsertAfter(" ")
$title.InsertAfter("ethical")
$title.InsertAfter(" ")
$title.InsertAfter("dilemmas")
$title.InsertAfter(" ")
$title.InsertAfter("you")
$title.InsertAfter(" ")
$title.InsertAfter("have")
$title.InsertAfter(" ")
$title.InsertAfter("faced")
$title.InsertAfter(" ")
$title.InsertAfter("in")
$title.InsertAfter(" ")
$title.InsertAfter("your")
$title.InsertAfter(" ")
$title.InsertAfter("job?")
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("Asking for the impossible")
$body.InsertAfter("`r")
$body.InsertAfter("“Blood from a turnip” test")
$body.Paragraphs(1).IndentLevel = 2
$body.Paragraphs(2).IndentLevel = 3

